$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")

# The "Destination" cell now points to the Results sheet cell (instead of Sheet2:C3)
$ws.Range("C7").Value = "Results:A1"

# Insert a new row for "Destination File" just below "Destination" and fill it in.
$ws.Rows.Item(8).Insert()
$ws.Range("B8").Value = "Destination File"
$ws.Range("C8").Value = "Results.xlsx"

# Keep the current selection where the user last clicked.
[void]$ws.Range("C8").Select()
